# Add a new column BB (quarterly date 2025-11-25 / serial 45986) that is a
# copy of column BA's QoQ series for the existing rows, with its own new
# value at row 83, plus a brand-new row 84 (serial 45884) that only has a
# value in BB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: BB1 = new quarter-end date, same date style as BA1 ---
$ws.Range("BA1").Copy()
$ws.Range("BB1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("BB1").Value = 45986

# --- Body rows 2-82: BB is a straight copy of BA's values (no style) ---
$ws.Range("BA2:BA82").Copy()
$ws.Range("BB2:BB82").PasteSpecial(-4163)   # xlPasteValues

# --- Row 83: BB83 gets its own (different) value, no style ---
$ws.Range("BB83").Value = -0.5

# --- New row 84: only A84 (date) and BB84 (value) are populated ---
$ws.Range("A83").Copy()
$ws.Range("A84").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A84").Value = 45884
$ws.Range("BB84").Value = -0.3

$excel.CutCopyMode = 0
